$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2020 column values (column Q)
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 4.4631700362051845
$ws.Range("Q6").Value = 22.107243650047039
$ws.Range("Q7").Value = 4.8469387755102042
$ws.Range("Q8").Value = 11.270912826533607
$ws.Range("Q9").Value = 8.2663605051664764
$ws.Range("Q10").Value = 9.0160381447767666
$ws.Range("Q11").Value = 2.7624309392265194
$ws.Range("Q12").Value = 1.1408815903197926
$ws.Range("Q13").Value = 1.7541111981205952
$ws.Range("Q14").Value = 3.6288232244686367

# Match formatting of the preceding (2019 / column P) cells
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)

$ws.Range("P6:P13").Copy()
$ws.Range("Q6:Q13").PasteSpecial(-4122)

$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Selection moves to P15 as in the saved file
$ws.Range("P15").Select() | Out-Null
